$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Counts (cps)" (column B) and "Error (cps)" (column C) values
# for each sample row (rows 2-15), reflecting the updated SCM results.

$data = @{
    2  = @{ B = 20.30291666666667;  C = 0.3677365781250001 }
    3  = @{ B = 20.27395833333333;  C = 0.36797234375 }
    4  = @{ B = 19.29708333333333;  C = 0.35892575 }
    5  = @{ B = 19.72104166666666;  C = 0.362374140625 }
    6  = @{ B = 20.00041666666667;  C = 0.365257609375 }
    7  = @{ B = 15.34416666666667;  C = 0.319925875 }
    8  = @{ B = 15.99166666666667;  C = 0.3264298958333333 }
    9  = @{ B = 15.614375;          C = 0.322827203125 }
    10 = @{ B = 13.18;              C = 0.29671475 }
    11 = @{ B = 12.936875;          C = 0.2935053515625 }
    12 = @{ B = 13.50416666666667;  C = 0.2999613020833333 }
    13 = @{ B = 17.84479166666667;  C = 0.34529671875 }
    14 = @{ B = 17.93;              C = 0.34560075 }
    15 = @{ B = 17.56416666666667;  C = 0.3422816979166666 }
}

foreach ($row in $data.Keys) {
    $ws.Range("B$row").Value = $data[$row].B
    $ws.Range("C$row").Value = $data[$row].C
}

$wb.Save()
